# Update latest output (run 106)
$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet: Cost ($) and Unit Cost ($/ML) for rows 3-5 ---
$wsSchedule.Range("E3").Value = 52.94103750000002
$wsSchedule.Range("F3").Value = 1.474270050125314

$wsSchedule.Range("E4").Value = 456.8181149999999
$wsSchedule.Range("F4").Value = 30.21283829365079

$wsSchedule.Range("E5").Value = -61.9336575
$wsSchedule.Range("F5").Value = -1.820507275132275

# --- Detailed sheet: Price (col B) and Type (col C) updates ---
$wsDetailed.Range("B33").Value = 0.7
$wsDetailed.Range("B34").Value = 0.80565

$wsDetailed.Range("B35").Value = -5.16887
$wsDetailed.Range("C35").Value = "historical"

$wsDetailed.Range("B36").Value = -3.42323
$wsDetailed.Range("C36").Value = "historical"

$wsDetailed.Range("B37").Value = -0.99633
$wsDetailed.Range("C37").Value = "historical"

$wsDetailed.Range("B38").Value = 9.87266
$wsDetailed.Range("B39").Value = 38.23834
$wsDetailed.Range("B40").Value = 59.50248
$wsDetailed.Range("B41").Value = 63.57365
$wsDetailed.Range("B42").Value = 64.46451
$wsDetailed.Range("B43").Value = 65
$wsDetailed.Range("B44").Value = 67.39644
$wsDetailed.Range("B45").Value = 65.14413999999999

$wsDetailed.Range("B47").Value = 58.35781
$wsDetailed.Range("B48").Value = 58.02218
$wsDetailed.Range("B49").Value = 63.67826
$wsDetailed.Range("B50").Value = 61.53081
$wsDetailed.Range("B51").Value = 57.06003
$wsDetailed.Range("B52").Value = 57.06003

$wsDetailed.Range("B59").Value = 57.13092

$wsDetailed.Range("B61").Value = 67.7428
$wsDetailed.Range("B62").Value = 67.60051
$wsDetailed.Range("B63").Value = 67.13655

$wsDetailed.Range("B65").Value = 8.78514
$wsDetailed.Range("B66").Value = 0.7
$wsDetailed.Range("B67").Value = 0.01001
$wsDetailed.Range("B68").Value = 0.00953
$wsDetailed.Range("B69").Value = -5.50985
$wsDetailed.Range("B70").Value = -6.58192
$wsDetailed.Range("B71").Value = -7.59468
$wsDetailed.Range("B72").Value = -9.99
$wsDetailed.Range("B73").Value = -8.222060000000001
$wsDetailed.Range("B74").Value = -9.70499
$wsDetailed.Range("B75").Value = -10.47783
$wsDetailed.Range("B76").Value = -9.99
$wsDetailed.Range("B77").Value = -9.55391
$wsDetailed.Range("B78").Value = -8.21705
$wsDetailed.Range("B79").Value = -9.465059999999999
$wsDetailed.Range("B80").Value = -7.40455
$wsDetailed.Range("B81").Value = -6.37448
$wsDetailed.Range("B82").Value = -1.25798
$wsDetailed.Range("B83").Value = -2.72939
$wsDetailed.Range("B84").Value = -0.1267
$wsDetailed.Range("B85").Value = 9.53589
$wsDetailed.Range("B86").Value = 9.66465
$wsDetailed.Range("B87").Value = 54.14803
$wsDetailed.Range("B88").Value = 64.8901
$wsDetailed.Range("B89").Value = 64.8901
$wsDetailed.Range("B90").Value = 64.8901
$wsDetailed.Range("B91").Value = 59.32108
$wsDetailed.Range("B92").Value = 58.56705
$wsDetailed.Range("B93").Value = 57.98785
$wsDetailed.Range("B94").Value = 59.44475
$wsDetailed.Range("B95").Value = 60.44065
$wsDetailed.Range("B96").Value = 62.50152
$wsDetailed.Range("B97").Value = 57.06003
